$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update quarterly income-statement figures (rows 11-27, columns D:M)
# Values sourced from the updated database export.

# Row 11
$ws.Range("D11").Value = 4211766
$ws.Range("E11").Value = 4770866
$ws.Range("F11").Value = 4264798
$ws.Range("G11").Value = 4584454
$ws.Range("H11").Value = 6698553
$ws.Range("I11").Value = 8843755
$ws.Range("J11").Value = 9477290
$ws.Range("K11").Value = 6565990
$ws.Range("L11").Value = 9090948
$ws.Range("M11").Value = 8495313

# Row 12
$ws.Range("D12").Value = -2914191
$ws.Range("E12").Value = -3047037
$ws.Range("F12").Value = -3274581
$ws.Range("G12").Value = -3256014
$ws.Range("H12").Value = -5377109
$ws.Range("I12").Value = -7520488
$ws.Range("J12").Value = -8070666
$ws.Range("K12").Value = -5652744
$ws.Range("L12").Value = -8807341
$ws.Range("M12").Value = -8023441

# Row 13
$ws.Range("D13").Value = 1297575
$ws.Range("E13").Value = 1723829
$ws.Range("F13").Value = 990217
$ws.Range("G13").Value = 1328440
$ws.Range("H13").Value = 1321444
$ws.Range("I13").Value = 1323267
$ws.Range("J13").Value = 1406624
$ws.Range("K13").Value = 913246
$ws.Range("L13").Value = 283607
$ws.Range("M13").Value = 471872

# Row 14
$ws.Range("D14").Value = 20687
$ws.Range("E14").Value = -43655
$ws.Range("F14").Value = -70928
$ws.Range("G14").Value = -64524
$ws.Range("H14").Value = -132440
$ws.Range("I14").Value = -81004
$ws.Range("J14").Value = -104145
$ws.Range("K14").Value = -135186
$ws.Range("L14").Value = -145510
$ws.Range("M14").Value = -113676

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = -2219
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("D16").Value = 124678
$ws.Range("E16").Value = 73481
$ws.Range("F16").Value = 67655
$ws.Range("G16").Value = 25557
$ws.Range("H16").Value = 95312
$ws.Range("I16").Value = 72026
$ws.Range("J16").Value = 142454
$ws.Range("K16").Value = 117730
$ws.Range("L16").Value = 77293
$ws.Range("M16").Value = 117461

# Row 17
$ws.Range("D17").Value = 1442940
$ws.Range("E17").Value = 1753655
$ws.Range("F17").Value = 986944
$ws.Range("G17").Value = 1289473
$ws.Range("H17").Value = 1284316
$ws.Range("I17").Value = 1314289
$ws.Range("J17").Value = 1442714
$ws.Range("K17").Value = 895790
$ws.Range("L17").Value = 215390
$ws.Range("M17").Value = 475657

# Row 18
$ws.Range("D18").Value = -39
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = -42305
$ws.Range("H18").Value = -56578
$ws.Range("I18").Value = -134264
$ws.Range("J18").Value = -157214
$ws.Range("K18").Value = -107755
$ws.Range("L18").Value = -117006
$ws.Range("M18").Value = -138699

# Row 19
$ws.Range("D19").Value = -14108
$ws.Range("E19").Value = 103948
$ws.Range("F19").Value = 70448
$ws.Range("G19").Value = 4902
$ws.Range("H19").Value = 97861
$ws.Range("I19").Value = 20656
$ws.Range("J19").Value = -43329
$ws.Range("K19").Value = 70201
$ws.Range("L19").Value = 102009
$ws.Range("M19").Value = 260392

# Row 20
$ws.Range("D20").Value = 1428793
$ws.Range("E20").Value = 1857603
$ws.Range("F20").Value = 1057392
$ws.Range("G20").Value = 1252070
$ws.Range("H20").Value = 1325599
$ws.Range("I20").Value = 1200681
$ws.Range("J20").Value = 1242171
$ws.Range("K20").Value = 858236
$ws.Range("L20").Value = 200393
$ws.Range("M20").Value = 597350

# Row 21
$ws.Range("D21").Value = -291243
$ws.Range("E21").Value = -414501
$ws.Range("F21").Value = 335943
$ws.Range("G21").Value = -273465
$ws.Range("H21").Value = -264824
$ws.Range("I21").Value = -237994
$ws.Range("J21").Value = 115412
$ws.Range("K21").Value = -179749
$ws.Range("L21").Value = 30165
$ws.Range("M21").Value = -89324

# Row 22
$ws.Range("D22").Value = 1137550
$ws.Range("E22").Value = 1443102
$ws.Range("F22").Value = 1393335
$ws.Range("G22").Value = 978605
$ws.Range("H22").Value = 1060775
$ws.Range("I22").Value = 962687
$ws.Range("J22").Value = 1357583
$ws.Range("K22").Value = 678487
$ws.Range("L22").Value = 230558
$ws.Range("M22").Value = 508026

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("D24").Value = 1137550
$ws.Range("E24").Value = 1443102
$ws.Range("F24").Value = 1393335
$ws.Range("G24").Value = 978605
$ws.Range("H24").Value = 1060775
$ws.Range("I24").Value = 962687
$ws.Range("J24").Value = 1357583
$ws.Range("K24").Value = 678487
$ws.Range("L24").Value = 230558
$ws.Range("M24").Value = 508026

# Row 25
$ws.Range("D25").Value = 875
$ws.Range("E25").Value = 1110
$ws.Range("F25").Value = 1072
$ws.Range("G25").Value = 753
$ws.Range("H25").Value = 816
$ws.Range("I25").Value = 741
$ws.Range("J25").Value = 1044
$ws.Range("K25").Value = 522
$ws.Range("L25").Value = 177
$ws.Range("M25").Value = 130

# Row 26
$ws.Range("D26").Value = 1300000
$ws.Range("E26").Value = 1300000
$ws.Range("F26").Value = 1300000
$ws.Range("G26").Value = 1300000
$ws.Range("H26").Value = 1300000
$ws.Range("I26").Value = 1300000
$ws.Range("J26").Value = 1300000
$ws.Range("K26").Value = 1300000
$ws.Range("L26").Value = 1300000
$ws.Range("M26").Value = 3900000

# Row 27
$ws.Range("D27").Value = 292
$ws.Range("E27").Value = 370
$ws.Range("F27").Value = 357
$ws.Range("G27").Value = 251
$ws.Range("H27").Value = 272
$ws.Range("I27").Value = 247
$ws.Range("J27").Value = 348
$ws.Range("K27").Value = 174
$ws.Range("L27").Value = 59
$ws.Range("M27").Value = 130
